$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old column N (which held the "Лаб5" data).
# This shifts the old N column (header "Лаб5" + its data cells) to column P,
# preserving values and styles, and creates blank N/O columns inheriting
# the per-row formatting (so header row cells land with the bold header
# style, data rows with whatever fill/border the row already used).
$ws.Columns("N:O").Insert()

# The Insert() leaves every row in N:O as an empty-but-styled cell. The
# target workbook only has real cells in N2/O2 (new headers) and N4/O4
# (new data), so clear the rest of the inserted range back to blank.
$ws.Range("N3:O31").Clear()

# New header cells for the two newly inserted lab columns.
$ws.Range("N2").Value = "Лаб3"
$ws.Range("O2").Value = "Лаб4"

# New grade entries in row 4 (student "Адеев Григорий").
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 5

# New grade entries elsewhere in the sheet.
$ws.Range("K6").Value = 5
$ws.Range("P6").Value = 5

$ws.Range("M13").Value = 5

$ws.Range("K22").Value = 5

$ws.Range("M23").Value = 5

$ws.Range("L25").Value = 5

$ws.Range("K30").Value = 5

# Move the active-cell selection to P6 (matches the saved view state).
[void]$ws.Range("P6").Select()
